$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# 1) Re-size the table's grid columns (w:tblGrid/w:gridCol widths, in
#    twentieths-of-a-point i.e. twips -> points for the COM API).
# ---------------------------------------------------------------------
$newWidthsTwips = @(1213, 3084, 618, 655, 770, 1284, 1442)
for ($i = 1; $i -le $t.Columns.Count; $i++) {
    $t.Columns.Item($i).Width = $newWidthsTwips[$i - 1] / 20.0
}

# ---------------------------------------------------------------------
# 2) The header row's numeric-column headings (Min / Max / Mean /
#    Std. Deviation / Ireland (mean)) switch from right- to
#    center-justified.
# ---------------------------------------------------------------------
$headerRow = $t.Rows.Item(1)
for ($c = 1; $c -le $headerRow.Cells.Count; $c++) {
    $cell = $headerRow.Cells.Item($c)
    $para = $cell.Range.Paragraphs.Item(1)
    if ($para.Alignment -eq 2) {
        $para.Alignment = 1
    }
}

# ---------------------------------------------------------------------
# 3) "Ireland<linebreak>                    (mean)" loses the manual
#    line break and two of its leading spaces:
#    "Ireland                  (mean)".
# ---------------------------------------------------------------------
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Rows.Item($r).Cells.Count; $c++) {
        $cell = $t.Rows.Item($r).Cells.Item($c)
        if ($cell.Range.Text.Contains("(mean)")) {
            $lineBreak = [char]11
            $old = $lineBreak + "                    (mean)"
            $new = "                  (mean)"
            $cell.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
        }
    }
}

# ---------------------------------------------------------------------
# 4 & 5) Minor sample-name corrections.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Laois County Council", $true, $false, $false, $false, $false, $true, 1, $false, "Laois County Council", 2) | Out-Null
$d.Content.Find.Execute("Post leaving cert", $true, $false, $false, $false, $false, $true, 1, $false, "Post leaving cert", 2) | Out-Null

Write-Host "edits applied"
